# FEAT: adds support for scanning base 16 binary! values.
#
# Adds a new "bin16 transitions" worksheet (after "transitions"), populates
# its FSM transition table, introduces a new bold/red "error" cell style,
# and makes a couple of small corrections to the existing "transitions"
# sheet (J12/K12 swap, Q13 fix).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Add the new worksheet right after "transitions" and make it the active
#    (selected) tab, matching the authored workbook.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "bin16 transitions"

# ---------------------------------------------------------------------------
# 2. Seed the new shared strings in the exact order they are first used so
#    the shared string table lines up with the authored file (101..110).
# ---------------------------------------------------------------------------
$ws2.Range("B1").Value = "C_BIN_SKIP"     # 101
$ws2.Range("C1").Value = "C_BIN_BLANK"    # 102
$ws2.Range("D1").Value = "C_BIN_LINE"     # 103
$ws2.Range("E1").Value = "C_BIN_HEXA"     # 104
$ws2.Range("E2").Value = "S_BIN_1ST"      # 105
$ws2.Range("A2").Value = "S_BIN_START"    # 106
$ws2.Range("E3").Value = "T_BIN_BYTE"     # 107
$ws2.Range("B3").Value = "T_BIN_ERROR"    # 108
$ws2.Range("F2").Value = "S_BIN_CMT"      # 109
$ws2.Range("F1").Value = "C_BIN_CMT"      # 110

# ---------------------------------------------------------------------------
# 3. Fill in the rest of the transition grid (reusing the strings above).
# ---------------------------------------------------------------------------
$ws2.Range("B2").Value = "S_BIN_START"
$ws2.Range("C2").Value = "S_BIN_START"
$ws2.Range("D2").Value = "S_BIN_START"

$ws2.Range("A3").Value = "S_BIN_1ST"
$ws2.Range("C3").Value = "T_BIN_ERROR"
$ws2.Range("D3").Value = "T_BIN_ERROR"
$ws2.Range("F3").Value = "T_BIN_ERROR"

$ws2.Range("A4").Value = "S_BIN_CMT"
$ws2.Range("B4").Value = "S_BIN_CMT"
$ws2.Range("C4").Value = "S_BIN_CMT"
$ws2.Range("D4").Value = "S_BIN_START"
$ws2.Range("E4").Value = "S_BIN_CMT"
$ws2.Range("F4").Value = "S_BIN_CMT"

# Row 5 / A5 is present but left blank (style only, set further down).

# ---------------------------------------------------------------------------
# 4. Apply formatting. Re-use the existing cell styles from "transitions" by
#    copy/pasting formats so we don't fork near-duplicate style entries.
# ---------------------------------------------------------------------------

# Header row (border + bold) -- matches sheet1 A1 / B1 styles.
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122)
$ws1.Range("B1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)

# Column A "state" cells -- matches sheet1 A2 style.
$ws1.Range("A2").Copy()
$ws2.Range("A2:A5").PasteSpecial(-4122)

# Plain sz8 cells -- matches sheet1 B2 style.
$ws1.Range("B2").Copy()
$ws2.Range("B2:F2").PasteSpecial(-4122)
$ws2.Range("B4:F4").PasteSpecial(-4122)

# Bold sz8 cell (T_BIN_BYTE) -- matches sheet1 F2 style.
$ws1.Range("F2").Copy()
$ws2.Range("E3").PasteSpecial(-4122)

# New bold / red "error" style, built once on B3 then copied onto the other
# T_BIN_ERROR cells (this ordering keeps the style table minimal: it adds
# exactly one new font + one new cell format, matching the authored file).
$ws2.Range("B3").HorizontalAlignment = -4131
$ws2.Range("B3").Font.Size = 8
$ws2.Range("B3").Font.Bold = $true
$ws2.Range("B3").Font.Color = 255
$ws2.Range("B3").Copy()
$ws2.Range("C3:D3").PasteSpecial(-4122)
$ws2.Range("F3").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5. View state: new sheet is the active/selected tab with G11 selected;
#    "transitions" loses its tabSelected flag and its selection moves to B2.
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("B2").Select()

$ws2.Activate()
$ws2.Range("G11").Select()

# ---------------------------------------------------------------------------
# 6. Small data corrections on "transitions": J12/K12 swap their contents,
#    and Q13 is corrected from C_WORD(65) to T_ERROR(8).
# ---------------------------------------------------------------------------
$ws1.Range("J12").Value = "C_BIN"
$ws1.Range("K12").Value = "C_WORD"
$ws1.Range("Q13").Value = "T_ERROR"
